# "updated boot for SSID and PSTR"
#
# 1) Table S1 - Plasticity AIC: refresh bootstrap-derived summary stats
#    (AIC/BIC/R2/ICC/RMSE/Sigma columns) for the SSID (pstr-model) rows 8-12.
# 2) Table S3 - Plasticity GLM: the PSTR coefficient block had a duplicated
#    "pCO2" row; drop the spurious row and refresh the remaining PSTR
#    coefficients to the corrected bootstrap values (PAST block is untouched
#    and simply shifts up).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Table S1 - Plasticity AIC
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table S1 - Plasticity AIC")

$ws1.Cells.Item(8, 4).Value  = 105.6
$ws1.Cells.Item(8, 5).Value  = 121.5
$ws1.Cells.Item(8, 6).Value  = 0.3969
$ws1.Cells.Item(8, 7).Value  = 0.317
$ws1.Cells.Item(8, 8).Value  = 0.117033435887396
$ws1.Cells.Item(8, 9).Value  = 1.02875627247403
$ws1.Cells.Item(8, 10).Value = 0.342860527251939

$ws1.Cells.Item(9, 4).Value  = 99.8
$ws1.Cells.Item(9, 5).Value  = 109.6
$ws1.Cells.Item(9, 6).Value  = 0.2918
$ws1.Cells.Item(9, 7).Value  = 0.2233
$ws1.Cells.Item(9, 8).Value  = 0.0882187502302896
$ws1.Cells.Item(9, 9).Value  = 1.15145387403649
$ws1.Cells.Item(9, 10).Value = 0.362141524711882
$ws1.Cells.Item(9, 11).Value = 47.1

$ws1.Cells.Item(10, 4).Value  = 100.8
$ws1.Cells.Item(10, 5).Value  = 111.8
$ws1.Cells.Item(10, 6).Value  = 0.3093
$ws1.Cells.Item(10, 7).Value  = 0.2614
$ws1.Cells.Item(10, 8).Value  = 0.0648056876434838
$ws1.Cells.Item(10, 9).Value  = 1.14657450978296
$ws1.Cells.Item(10, 10).Value = 0.36062662799942
$ws1.Cells.Item(10, 11).Value = 45

$ws1.Cells.Item(11, 4).Value  = 101.6
$ws1.Cells.Item(11, 5).Value  = 112.6
$ws1.Cells.Item(11, 6).Value  = 0.2784
$ws1.Cells.Item(11, 7).Value  = 0.2378
$ws1.Cells.Item(11, 8).Value  = 0.0533714820211562
$ws1.Cells.Item(11, 9).Value  = 1.13665595795787
$ws1.Cells.Item(11, 10).Value = 0.359834556845443
$ws1.Cells.Item(11, 11).Value = 36.2

$ws1.Cells.Item(12, 4).Value  = 97.5
$ws1.Cells.Item(12, 5).Value  = 104.8
$ws1.Cells.Item(12, 6).Value  = 0.232
$ws1.Cells.Item(12, 7).Value  = 0.188
$ws1.Cells.Item(12, 8).Value  = 0.0541808464958402
$ws1.Cells.Item(12, 9).Value  = 1.22142050234447
$ws1.Cells.Item(12, 10).Value = 0.371319943278335
$ws1.Cells.Item(12, 11).Value = 28.8

# ---------------------------------------------------------------------------
# Table S3 - Plasticity GLM
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table S3 - Plasticity GLM")

# Old row 12 (A="pCO2", B=0.194, C=0.418, D=0.46, E=0.643) was a spurious
# duplicate entry in the PSTR block; remove it and let rows 13-20 shift up.
$ws3.Rows("12:12").Delete()

# Remaining PSTR rows (now 11-14) get refreshed coefficient estimates.
$ws3.Cells.Item(11, 3).Value = 0.148
$ws3.Cells.Item(11, 4).Value = 8.66

$ws3.Cells.Item(12, 2).Value = -0.338
$ws3.Cells.Item(12, 3).Value = 0.193
$ws3.Cells.Item(12, 4).Value = -1.75
$ws3.Cells.Item(12, 5).Value = 0.08

$ws3.Cells.Item(13, 2).Value = -0.059
$ws3.Cells.Item(13, 3).Value = 0.187
$ws3.Cells.Item(13, 4).Value = -0.31
$ws3.Cells.Item(13, 5).Value = 0.753

$ws3.Cells.Item(14, 2).Value = 0.227
$ws3.Cells.Item(14, 3).Value = 0.173
$ws3.Cells.Item(14, 4).Value = 1.31
$ws3.Cells.Item(14, 5).Value = 0.19
